# Update "想去人数" (want-to-go count) values in column F across sheets.
# Values correspond to a refreshed data pull (gh-pages output regenerated).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1168
$ws1.Range("F3").Value = 1081
$ws1.Range("F4").Value = 1884
$ws1.Range("F10").Value = 324
$ws1.Range("F13").Value = 776
$ws1.Range("F14").Value = 215
$ws1.Range("F17").Value = 116
$ws1.Range("F19").Value = 196
$ws1.Range("F25").Value = 897
$ws1.Range("F26").Value = 340
$ws1.Range("F29").Value = 297
$ws1.Range("F32").Value = 418

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F11").Value = 127

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1168
$ws4.Range("F4").Value = 1081
$ws4.Range("F5").Value = 1884
$ws4.Range("F12").Value = 324
$ws4.Range("F15").Value = 776
$ws4.Range("F16").Value = 215
$ws4.Range("F21").Value = 116
$ws4.Range("F27").Value = 196
$ws4.Range("F33").Value = 897
$ws4.Range("F34").Value = 340
$ws4.Range("F39").Value = 297
$ws4.Range("F41").Value = 127
$ws4.Range("F42").Value = 127
$ws4.Range("F46").Value = 418
